$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: Foster's Home for Imaginary Friends ---
# --- Row 12: Rocko's Modern Life ---
# --- Row 13: Ducktales ---
# Shared strings are created in the Excel shared-string table in the order
# new unique values are first written to a cell, so we set the cells that
# introduce brand-new text in the same order the source workbook used.

# New titles first
$ws.Range("A11").Value = "Foster's Home for Imaginary Friends"
$ws.Range("A12").Value = "Rocko's Modern Life"
$ws.Range("A13").Value = "Ducktales"

# Row 11 remaining new text values
$ws.Range("F11").Value = "Fantasy, Comedy, Adventure"
$ws.Range("C11").Value = "A boy and his beloved imaginary friend are able to stay together at an orphanage of sorts for imaginary friends that children have outgrown to be adopted by new children."
$ws.Range("I11").Value = "https://m.media-amazon.com/images/M/MV5BNjYyNGFjOTctYzFmNC00NzdmLThhMDgtNjEzZTRmNzA3ODc5XkEyXkFqcGdeQXVyNjk1Njg5NTA@._V1_.jpg"

# Row 12 remaining new text values
$ws.Range("B12").Value = "Joe Murray"
$ws.Range("F12").Value = "Satire, Slapstick"
$ws.Range("C12").Value = "The wacky misadventures of an Australian wallaby and his friends as he finishes his transition to American life."
$ws.Range("I12").Value = "https://m.media-amazon.com/images/M/MV5BMjc5YmYwZmEtZjA3Ni00MWUxLWFjMmYtMzE3NTNiY2MyZTlmXkEyXkFqcGdeQXVyNjk1Njg5NTA@._V1_.jpg"

# Row 13 remaining new text values
$ws.Range("C13").Value = "The comedy-adventure series chronicles the high-flying adventures of trillionaire Scrooge McDuck, his grandnephews - Huey, Dewey, and Louie - temperamental nephew Donald Duck, Launchpad McQuack, Mrs. Beakley, and her granddaughter Webby."
$ws.Range("B13").Value = "Matt Youngberg, Francisco Angones"
$ws.Range("I13").Value = "https://m.media-amazon.com/images/M/MV5BNTA2NTc5MzQwNV5BMl5BanBnXkFtZTgwOTY2ODI2MjI@._V1_SY1000_CR0,0,666,1000_AL_.jpg"

# Remaining cells reuse existing shared strings / are plain numbers.
$ws.Range("B11").Value = "Craig McCracken"
$ws.Range("D11").Value = 2004
$ws.Range("E11").Value = 2009
$ws.Range("G11").Value = 6
$ws.Range("H11").Value = 79
$ws.Range("J11").Value = "Cartoon Network"

$ws.Range("D12").Value = 1993
$ws.Range("E12").Value = 1996
$ws.Range("G12").Value = 4
$ws.Range("H12").Value = 52
$ws.Range("J12").Value = "Nickelodeon"

$ws.Range("D13").Value = 2017
$ws.Range("F13").Value = "Comedy, Adventure, Science Fantasy"
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 43
$ws.Range("J13").Value = "Disney Channel, Disney XD"

# Hyperlinks for the new Image cells (match existing rows' hyperlink style)
$ws.Hyperlinks.Add($ws.Range("I11"), "https://m.media-amazon.com/images/M/MV5BNjYyNGFjOTctYzFmNC00NzdmLThhMDgtNjEzZTRmNzA3ODc5XkEyXkFqcGdeQXVyNjk1Njg5NTA@._V1_.jpg") | Out-Null
$ws.Range("I11").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("I12"), "https://m.media-amazon.com/images/M/MV5BMjc5YmYwZmEtZjA3Ni00MWUxLWFjMmYtMzE3NTNiY2MyZTlmXkEyXkFqcGdeQXVyNjk1Njg5NTA@._V1_.jpg") | Out-Null
$ws.Range("I12").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("I13"), "https://m.media-amazon.com/images/M/MV5BNTA2NTc5MzQwNV5BMl5BanBnXkFtZTgwOTY2ODI2MjI@._V1_SY1000_CR0,0,666,1000_AL_.jpg") | Out-Null
$ws.Range("I13").Style = "Hyperlink"

# Match the active selection shown in the target workbook
$ws.Range("A14").Select() | Out-Null
